$d = $word.ActiveDocument

# Locate the "Operating hours" paragraph (OH = Operating hours: ...) by its
# distinctive text so we don't depend on a brittle paragraph index.
$target = $null
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $cand = $paras.Item($i)
    if ($cand.Range.Text -like "*Operating hours:*weeks per year*") {
        $target = $cand
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the Operating hours paragraph"
}

$r = $target.Range

# Replace the whole paragraph (including its end-of-paragraph mark) with the
# updated OOXML: adds a hanging indent to the paragraph properties and
# rewrites the trailing "(...) hours per day, ... days per week, ... weeks
# per year)" text into abbreviated "hrs/day, .../wk, ...wks/yr" units spread
# across a tab-aligned layout, matching Word's normal run/proofErr markup.
$xml = '<w:p w14:paraId="32EC03CE" w14:textId="247346AE" w:rsidR="000E317C" w:rsidRPr="00892A30" w:rsidRDefault="000E317C" w:rsidP="00892A30"><w:pPr><w:pBdr><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/><w:between w:val="nil"/></w:pBdr><w:spacing w:line="360" w:lineRule="auto"/><w:ind w:left="720" w:hanging="720"/><w:jc w:val="both"/><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:tab/><w:t>OH</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve">= Operating hours: </w:t></w:r><w:r w:rsidR="00892A30"><w:t>${OH}</w:t></w:r><w:r w:rsidR="00892A30" w:rsidRPr="0027079E"><w:rPr><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00892A30" w:rsidRPr="0027079E"><w:rPr><w:color w:val="000000"/></w:rPr><w:t>hr</w:t></w:r><w:r w:rsidR="00892A30"><w:rPr><w:color w:val="000000"/></w:rPr><w:t>s</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00892A30" w:rsidRPr="0027079E"><w:rPr><w:color w:val="000000"/></w:rPr><w:t>/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00892A30" w:rsidRPr="0027079E"><w:rPr><w:color w:val="000000"/></w:rPr><w:t>yr</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00892A30"><w:rPr><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> (${HR} </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>hrs</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>/</w:t></w:r><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>day, ${DY} days</w:t></w:r><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>wk</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>, ${WK</w:t></w:r><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>}</w:t></w:r><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:tab/><w:t xml:space="preserve">    </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>wks</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>y</w:t></w:r><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>r</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>)</w:t></w:r></w:p>'

$r.InsertXML($xml)

Write-Host "Operating hours paragraph updated."
